$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 169.915657
$ws.Range("H2").Value = 509.746971
$ws.Range("I2").Value = 0.4441184931734509
$ws.Range("J2").Value = 0.4441184931734509
$ws.Range("M2").Value = 2.590303333333333
$ws.Range("N2").Value = 7.770910000000001
$ws.Range("O2").Value = 0.4994568961971165
$ws.Range("P2").Value = 0.4994568961971165
$ws.Range("Q2").Value = 440.1330927126234
$ws.Range("R2").Value = 3961.197834413611
$ws.Range("S2").Value = 0.221818044144152
$ws.Range("T2").Value = 0.2218180441441521
# Row 3
$ws.Range("G3").Value = 169.915657
$ws.Range("H3").Value = 509.746971
$ws.Range("I3").Value = 0.4441184931734509
$ws.Range("J3").Value = 0.4441184931734509
$ws.Range("O3").Value = 0.3915901179531479
$ws.Range("P3").Value = 0.3915901179531478
$ws.Range("Q3").Value = 345.0783661267124
$ws.Range("R3").Value = 3105.705295140411
$ws.Range("S3").Value = 0.173912413126966
$ws.Range("T3").Value = 0.1739124131269659
# Row 4
$ws.Range("G4").Value = 169.915657
$ws.Range("H4").Value = 509.746971
$ws.Range("I4").Value = 0.4441184931734509
$ws.Range("J4").Value = 0.4441184931734509
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5084619999999999
$ws.Range("N4").Value = 1.525386
$ws.Range("O4").Value = 0.09804058431541925
$ws.Range("P4").Value = 0.09804058431541923
$ws.Range("Q4").Value = 86.39565478953398
$ws.Range("R4").Value = 777.5608931058059
$ws.Range("S4").Value = 0.04354163657600866
$ws.Range("T4").Value = 0.04354163657600866
# Row 5
$ws.Range("G5").Value = 169.915657
$ws.Range("H5").Value = 509.746971
$ws.Range("I5").Value = 0.4441184931734509
$ws.Range("J5").Value = 0.4441184931734509
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05659433333333333
$ws.Range("N5").Value = 0.169783
$ws.Range("O5").Value = 0.01091240153431645
$ws.Range("P5").Value = 0.01091240153431645
$ws.Range("Q5").Value = 9.616263330810334
$ws.Range("R5").Value = 86.546369977293
$ws.Range("S5").Value = 0.004846399326324274
$ws.Range("T5").Value = 0.004846399326324274
# Row 6
$ws.Range("I6").Value = 0.1787346690539575
$ws.Range("J6").Value = 0.1787346690539575
$ws.Range("M6").Value = 2.590303333333333
$ws.Range("N6").Value = 7.770910000000001
$ws.Range("O6").Value = 0.4994568961971165
$ws.Range("P6").Value = 0.4994568961971165
$ws.Range("Q6").Value = 177.13075198371
$ws.Range("R6").Value = 1594.17676785339
$ws.Range("S6").Value = 0.08927026304850844
$ws.Range("T6").Value = 0.08927026304850842
# Row 7
$ws.Range("I7").Value = 0.1787346690539575
$ws.Range("J7").Value = 0.1787346690539575
$ws.Range("O7").Value = 0.3915901179531479
$ws.Range("P7").Value = 0.3915901179531478
$ws.Range("S7").Value = 0.06999073013715608
$ws.Range("T7").Value = 0.06999073013715607
# Row 8
$ws.Range("I8").Value = 0.1787346690539575
$ws.Range("J8").Value = 0.1787346690539575
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5084619999999999
$ws.Range("N8").Value = 1.525386
$ws.Range("O8").Value = 0.09804058431541925
$ws.Range("P8").Value = 0.09804058431541923
$ws.Range("Q8").Value = 34.76977204026599
$ws.Range("R8").Value = 312.9279483623939
$ws.Range("S8").Value = 0.01752325139147308
$ws.Range("T8").Value = 0.01752325139147308
# Row 9
$ws.Range("I9").Value = 0.1787346690539575
$ws.Range("J9").Value = 0.1787346690539575
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05659433333333333
$ws.Range("N9").Value = 0.169783
$ws.Range("O9").Value = 0.01091240153431645
$ws.Range("P9").Value = 0.01091240153431645
$ws.Range("Q9").Value = 3.870047454423
$ws.Range("R9").Value = 34.830427089807
$ws.Range("S9").Value = 0.001950424476819948
$ws.Range("T9").Value = 0.001950424476819948
# Row 10
$ws.Range("G10").Value = 53.27463399999999
$ws.Range("H10").Value = 159.823902
$ws.Range("I10").Value = 0.1392470275793777
$ws.Range("J10").Value = 0.1392470275793778
$ws.Range("M10").Value = 2.590303333333333
$ws.Range("N10").Value = 7.770910000000001
$ws.Range("O10").Value = 0.4994568961971165
$ws.Range("P10").Value = 0.4994568961971165
$ws.Range("Q10").Value = 137.9974620323133
$ws.Range("R10").Value = 1241.97715829082
$ws.Range("S10").Value = 0.06954788819947029
$ws.Range("T10").Value = 0.06954788819947029
# Row 11
$ws.Range("G11").Value = 53.27463399999999
$ws.Range("H11").Value = 159.823902
$ws.Range("I11").Value = 0.1392470275793777
$ws.Range("J11").Value = 0.1392470275793778
$ws.Range("O11").Value = 0.3915901179531479
$ws.Range("P11").Value = 0.3915901179531478
$ws.Range("Q11").Value = 108.1944064561313
$ws.Range("R11").Value = 973.7496581051819
$ws.Range("S11").Value = 0.05452775995443377
$ws.Range("T11").Value = 0.05452775995443376
# Row 12
$ws.Range("G12").Value = 53.27463399999999
$ws.Range("H12").Value = 159.823902
$ws.Range("I12").Value = 0.1392470275793777
$ws.Range("J12").Value = 0.1392470275793778
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.5084619999999999
$ws.Range("N12").Value = 1.525386
$ws.Range("O12").Value = 0.09804058431541925
$ws.Range("P12").Value = 0.09804058431541923
$ws.Range("Q12").Value = 27.08812695290799
$ws.Range("R12").Value = 243.7931425761719
$ws.Range("S12").Value = 0.01365185994806749
$ws.Range("T12").Value = 0.01365185994806749
# Row 13
$ws.Range("G13").Value = 53.27463399999999
$ws.Range("H13").Value = 159.823902
$ws.Range("I13").Value = 0.1392470275793777
$ws.Range("J13").Value = 0.1392470275793778
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.05659433333333333
$ws.Range("N13").Value = 0.169783
$ws.Range("O13").Value = 0.01091240153431645
$ws.Range("P13").Value = 0.01091240153431645
$ws.Range("Q13").Value = 3.015042394807333
$ws.Range("R13").Value = 27.13538155326599
$ws.Range("S13").Value = 0.001519519477406206
$ws.Range("T13").Value = 0.001519519477406206
# Row 14
$ws.Range("G14").Value = 91.01828266666666
$ws.Range("H14").Value = 273.054848
$ws.Range("I14").Value = 0.2378998101932138
$ws.Range("J14").Value = 0.2378998101932138
$ws.Range("M14").Value = 2.590303333333333
$ws.Range("N14").Value = 7.770910000000001
$ws.Range("O14").Value = 0.4994568961971165
$ws.Range("P14").Value = 0.4994568961971165
$ws.Range("Q14").Value = 235.7649609857422
$ws.Range("R14").Value = 2121.88464887168
$ws.Range("S14").Value = 0.1188207008049857
$ws.Range("T14").Value = 0.1188207008049857
# Row 15
$ws.Range("G15").Value = 91.01828266666666
$ws.Range("H15").Value = 273.054848
$ws.Range("I15").Value = 0.2378998101932138
$ws.Range("J15").Value = 0.2378998101932138
$ws.Range("O15").Value = 0.3915901179531479
$ws.Range("P15").Value = 0.3915901179531478
$ws.Range("Q15").Value = 184.8472402415076
$ws.Range("R15").Value = 1663.625162173568
$ws.Range("S15").Value = 0.09315921473459209
$ws.Range("T15").Value = 0.09315921473459207
# Row 16
$ws.Range("G16").Value = 91.01828266666666
$ws.Range("H16").Value = 273.054848
$ws.Range("I16").Value = 0.2378998101932138
$ws.Range("J16").Value = 0.2378998101932138
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.5084619999999999
$ws.Range("N16").Value = 1.525386
$ws.Range("O16").Value = 0.09804058431541925
$ws.Range("P16").Value = 0.09804058431541923
$ws.Range("Q16").Value = 46.27933804125865
$ws.Range("R16").Value = 416.5140423713279
$ws.Range("S16").Value = 0.02332383639987001
$ws.Range("T16").Value = 0.02332383639987001
# Row 17
$ws.Range("G17").Value = 91.01828266666666
$ws.Range("H17").Value = 273.054848
$ws.Range("I17").Value = 0.2378998101932138
$ws.Range("J17").Value = 0.2378998101932138
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.05659433333333333
$ws.Range("N17").Value = 0.169783
$ws.Range("O17").Value = 0.01091240153431645
$ws.Range("P17").Value = 0.01091240153431645
$ws.Range("Q17").Value = 5.151119028664889
$ws.Range("R17").Value = 46.360071257984
$ws.Range("S17").Value = 0.002596058253766018
$ws.Range("T17").Value = 0.002596058253766018

Write-Host "Applied all changes"